# Adds a new "Ready for handoff" row (row 3) to the Overview, zh-cn, and de-de
# sheets for the newly handed-off file, mirroring the existing row 2 entries.

$wb = $excel.ActiveWorkbook

$mdName      = '4d37533f-a4df-45e9-8235-59d9d62483a0ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$mdDisplay   = 'e2e\4d37533f-a4df-45e9-8235-59d9d62483a0ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$zhcnXlf     = '4d37533f-a4df-45e9-8235-59d9d62483a0oooooooooooooooooooooooooooooooooooooooo.3610d29cb80748a1fc544b5622e4b14b9d8dcc59.zh-cn.xlf'
$dedeXlf     = '4d37533f-a4df-45e9-8235-59d9d62483a0oooooooooooooooooooooooooooooooooooooooo.3610d29cb80748a1fc544b5622e4b14b9d8dcc59.de-de.xlf'
$hoDate      = '2016-08-27 18:35:41'
$handoffDateZh = '2016-08-27 18:35:37'
$hyperlinkUrl = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/35520ef7b383221aa0a32e3ef102ea155f2bc917/e2e/4d37533f-a4df-45e9-8235-59d9d62483a0ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$dtFormat    = 'yyyy-mm-dd HH:mm:ss'

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $hoDate
$wsOverview.Range("G3").NumberFormat = $dtFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkUrl, "", "", $mdDisplay)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

$wsOverview.Range("E1:E3").EntireColumn.AutoFit()
$wsOverview.Range("F1:F3").EntireColumn.AutoFit()

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $zhcnXlf
$wsZhCn.Range("H3").Value = $handoffDateZh
$wsZhCn.Range("H3").NumberFormat = $dtFormat
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = $dtFormat
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkUrl, "", "", $mdDisplay)

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

$wsZhCn.Range("C1:C3").EntireColumn.AutoFit()

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $dedeXlf
$wsDeDe.Range("H3").Value = $hoDate
$wsDeDe.Range("H3").NumberFormat = $dtFormat
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = $dtFormat
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkUrl, "", "", $mdDisplay)

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

$wsDeDe.Range("C1:C3").EntireColumn.AutoFit()

Write-Host "Handoff report row added."
